# Weekly price update: a new observation is inserted as the new row 36
# (Arveja Verde, Terminal Hortofrutícola Agro Chillán), pushing the former
# rows 36..127 down to 37..128.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 36; this shifts rows
# 36-127 down to 37-128 (matches the diff: old D36..D127 reappear as
# D37..D128, and the sheet dimension grows from A1:R127 to A1:R128).
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the new weekly data point.
$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value = "Ñuble"
$ws.Cells.Item(36, 4).Value = 45274
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = 100112022
$ws.Cells.Item(36, 7).Value = "Arveja Verde"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 100
$ws.Cells.Item(36, 11).Value = 22000
$ws.Cells.Item(36, 12).Value = 22000
$ws.Cells.Item(36, 13).Value = 22000
$ws.Cells.Item(36, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(36, 15).Value = "Región de Ñuble"
$ws.Cells.Item(36, 16).Value = 880
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of
# column D (style index 2 in the original workbook).
$ws.Cells.Item(36, 4).NumberFormat = $ws.Cells.Item(37, 4).NumberFormat
